$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed Price (D) and Volume(1h) (E) cells, forcing text storage
# so numeric-looking strings (e.g. "599.25") are not coerced into numbers,
# matching the original inlineStr text cells.
$cells = @{
    'D2' = '67.742.46'
    'E2' = '  +0.07%  '
    'D3' = '3.788.29'
    'E3' = '  -0.11%  '
    'E4' = '  -0.07%  '
    'D5' = '599.25'
    'E5' = '  +0.54%  '
    'D6' = '165.20'
    'E6' = '  -0.97%  '
    'D7' = '0.999'
    'E7' = '  -0.12%  '
    'D8' = '0.516'
    'E8' = '  -0.92%  '
    'E9' = '  -0.74%  '
    'D10' = '0.451'
    'E10' = '  +0.34%  '
    'E11' = '  +1.20%  '
    'D12' = '0.0000249'
    'E12' = '  -1.40%  '
    'D13' = '35.67'
    'E13' = '  -1.09%  '
    'D14' = '4.420.43'
    'E14' = '  -0.17%  '
    'D15' = '3.773.51'
    'E15' = '  +0.85%  '
    'D16' = '67.714.13'
    'E16' = '  -0.01%  '
    'D17' = '18.32'
    'E17' = '  -1.54%  '
    'E18' = '  +1.77%  '
    'E19' = '  -0.19%  '
    'D20' = '461.07'
    'E20' = '  +0.26%  '
    'E21' = '  -2.46%  '
    'E22' = '  -0.66%  '
    'E23' = '  -3.36%  '
    'D24' = '82.64'
    'E24' = '  -0.83%  '
    'E25' = '  -0.22%  '
    'E26' = '  -0.01%  '
    'E27' = '  -0.04%  '
    'D28' = '9.94'
    'E28' = '  -0.80%  '
    'D29' = '3.935.26'
    'E29' = '  -0.11%  '
    'E30' = '  +3.22%  '
    'D31' = '2.62'
    'E31' = '  -5.48%  '
    'E32' = '  -2.93%  '
    'D33' = '29.10'
    'E33' = '  -1.75%  '
    'E34' = '  +0.13%  '
    'D35' = '8.97'
    'E35' = '  -1.13%  '
    'D36' = '0.0991'
    'E36' = '  -0.96%  '
    'E37' = '  +0.61%  '
    'E38' = '  -1.97%  '
    'D39' = '5.77'
    'E39' = '  -0.21%  '
    'D40' = '0.985'
    'E40' = '  -1.01%  '
    'E41' = '  -0.01%  '
    'D43' = '47.44'
    'E43' = '  -1.51%  '
    'D44' = '43.36'
    'E44' = '  -1.17%  '
    'E45' = '  +0.25%  '
    'D46' = '151.62'
    'E46' = '  +0.79%  '
    'D47' = '8.33'
    'E48' = '  +1.75%  '
    'D49' = '392.74'
    'E49' = '  +0.85%  '
    'E50' = '  +0.17%  '
    'D51' = '1.34'
    'E51' = '  +5.49%  '
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
}
